$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.6875
$ws.Range("P2").Value = 0.0625
$ws.Range("S2").Value = 0.0625
$ws.Range("C3").Value = 0.1538461538461539
$ws.Range("P3").Value = 0.4615384615384616
$ws.Range("S3").Value = 0.3846153846153846
$ws.Range("P5").Value = 1
$ws.Range("J6").Value = 0.3333333333333333
$ws.Range("Q6").Value = 0.1333333333333333
$ws.Range("R6").Value = 0.1333333333333333
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.05555555555555555
$ws.Range("Q7").Value = 0.2222222222222222
$ws.Range("S7").Value = 0.5555555555555556
$ws.Range("B8").Value = 0.1
$ws.Range("J8").Value = 0.1
$ws.Range("O8").Value = 0.1
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.1
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.04347826086956522
$ws.Range("J9").Value = 0.1739130434782609
$ws.Range("Q9").Value = 0.2608695652173913
$ws.Range("S9").Value = 0.5217391304347826
$ws.Range("B10").Value = 0.09876543209876543
$ws.Range("D10").Value = 0.04938271604938271
$ws.Range("E10").Value = 0.01234567901234568
$ws.Range("F10").Value = 0.08641975308641975
$ws.Range("J10").Value = 0.1358024691358025
$ws.Range("O10").Value = 0.03703703703703703
$ws.Range("Q10").Value = 0.1358024691358025
$ws.Range("R10").Value = 0.08641975308641975
$ws.Range("S10").Value = 0.3580246913580247
$ws.Range("F11").Value = 0.03225806451612903
$ws.Range("G11").Value = 0.1612903225806452
$ws.Range("J11").Value = 0.06451612903225806
$ws.Range("K11").Value = 0.2258064516129032
$ws.Range("L11").Value = 0.4838709677419355
$ws.Range("S11").Value = 0.03225806451612903
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.2
$ws.Range("S12").Value = 0.1333333333333333
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.25
$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("I15").Value = 0.1333333333333333
$ws.Range("J15").Value = 0.2666666666666667
$ws.Range("M15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.4
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.08333333333333333
$ws.Range("F17").Value = 0.03703703703703703
$ws.Range("H17").Value = 0.07407407407407407
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.3703703703703703
$ws.Range("K17").Value = 0.1111111111111111
$ws.Range("M17").Value = 0.03703703703703703
$ws.Range("O17").Value = 0.03703703703703703
$ws.Range("S17").Value = 0.2222222222222222
$ws.Range("H18").Value = 0.09090909090909091
$ws.Range("I18").Value = 0.2727272727272727
$ws.Range("J18").Value = 0.2727272727272727
$ws.Range("S18").Value = 0.3636363636363636
$ws.Range("F19").Value = 0.03669724770642202
$ws.Range("H19").Value = 0.1192660550458716
$ws.Range("I19").Value = 0.1192660550458716
$ws.Range("J19").Value = 0.3119266055045872
$ws.Range("K19").Value = 0.1651376146788991
$ws.Range("M19").Value = 0.01834862385321101
$ws.Range("O19").Value = 0.06422018348623854
$ws.Range("S19").Value = 0.1651376146788991
